# feat: add 2022-Q4 data
#
# Inserts a new worksheet "2022-Q4" right after "总计" (before the existing
# "2022-Q3" sheet), populated with the new quarter's fund holdings, and
# updates the "总计" (summary) sheet with a new leading row plus renumbered
# index column. The other pre-existing sheets ("2022-Q3", "2022-Q2",
# "2021-Q4", "2020-Q4") are left untouched - they simply shift right by one
# tab position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet by cloning the "2022-Q3" sheet (so it
#    inherits the same sheetPr/sheetView/style scaffolding) immediately
#    before it, then overwrite its contents.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$insertBefore = $wb.Worksheets.Item(2)
$templateSheet.Copy($insertBefore)

$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q4"

# Clear out whatever the template sheet had beyond the header row.
$ws.Range("A2:H100").ClearContents()

# Header row (row 1) - same labels as the other quarterly sheets.
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows. Fund codes / percentages are stored as text (leading zeros,
# fixed decimal places must be preserved) - a leading apostrophe forces
# text entry the same way typing it into Excel would, and the Style reset
# afterwards clears the quote-prefix formatting so the cell ends up with
# plain (unstyled) text, matching the other quarterly sheets.

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'002345"
$ws.Range("C2").Value = "华夏高端制造灵活配置混合A"
$ws.Range("D2").Value = "'16.65"
$ws.Range("E2").Value = "'93.68"
$ws.Range("F2").Value = "'4.24"
$ws.Range("G2").Value = "'0.7060"
$ws.Range("H2").Value = 10

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'009697"
$ws.Range("C3").Value = "华夏成长精选6个月定期开放混合A"
$ws.Range("D3").Value = "'5.66"
$ws.Range("E3").Value = "'92.44"
$ws.Range("F3").Value = "'4.42"
$ws.Range("G3").Value = "'0.2502"
$ws.Range("H3").Value = 7

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'006010"
$ws.Range("C4").Value = "国融融银灵活配置混合C"
$ws.Range("D4").Value = "'2.12"
$ws.Range("E4").Value = "'65.24"
$ws.Range("F4").Value = "'4.68"
$ws.Range("G4").Value = "'0.0992"
$ws.Range("H4").Value = 3

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'009698"
$ws.Range("C5").Value = "华夏成长精选6个月定期开放混合C"
$ws.Range("D5").Value = "'1.97"
$ws.Range("E5").Value = "'92.44"
$ws.Range("F5").Value = "'4.42"
$ws.Range("G5").Value = "'0.0871"
$ws.Range("H5").Value = 7

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'015058"
$ws.Range("C6").Value = "华夏高端制造灵活配置混合C"
$ws.Range("D6").Value = "'0.68"
$ws.Range("E6").Value = "'93.68"
$ws.Range("F6").Value = "'4.24"
$ws.Range("G6").Value = "'0.0288"
$ws.Range("H6").Value = 10

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'015429"
$ws.Range("C7").Value = "中银证券专精特新股票A"
$ws.Range("D7").Value = "'1.18"
$ws.Range("E7").Value = "'61.04"
$ws.Range("F7").Value = "'2.13"
$ws.Range("G7").Value = "'0.0251"
$ws.Range("H7").Value = 6

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'006009"
$ws.Range("C8").Value = "国融融银灵活配置混合A"
$ws.Range("D8").Value = "'0.19"
$ws.Range("E8").Value = "'65.24"
$ws.Range("F8").Value = "'4.68"
$ws.Range("G8").Value = "'0.0089"
$ws.Range("H8").Value = 3

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'015430"
$ws.Range("C9").Value = "中银证券专精特新股票C"
$ws.Range("D9").Value = "'0.14"
$ws.Range("E9").Value = "'61.04"
$ws.Range("F9").Value = "'2.13"
$ws.Range("G9").Value = "'0.0030"
$ws.Range("H9").Value = 6

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'005538"
$ws.Range("C10").Value = "中航新起航灵活配置混合C"
$ws.Range("D10").Value = "'0.04"
$ws.Range("E10").Value = "'69.18"
$ws.Range("F10").Value = "'6.47"
$ws.Range("G10").Value = "'0.0026"
$ws.Range("H10").Value = 3

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'005053"
$ws.Range("C11").Value = "银河量化价值混合A"
$ws.Range("D11").Value = "'0.11"
$ws.Range("E11").Value = "'70.26"
$ws.Range("F11").Value = "'1.78"
$ws.Range("G11").Value = "'0.0020"
$ws.Range("H11").Value = 6

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'005537"
$ws.Range("C12").Value = "中航新起航灵活配置混合A"
$ws.Range("D12").Value = "'0.02"
$ws.Range("E12").Value = "'69.18"
$ws.Range("F12").Value = "'6.47"
$ws.Range("G12").Value = "'0.0013"
$ws.Range("H12").Value = 3

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'013026"
$ws.Range("C13").Value = "银河量化价值混合C"
$ws.Range("D13").Value = "'0.00"
$ws.Range("E13").Value = "'70.26"
$ws.Range("F13").Value = "'1.78"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 6

# Drop the quote-prefix styling picked up from the leading apostrophes so
# the text-but-numeric-looking cells end up unstyled, like the source data.
$ws.Range("B2:G13").Style = "Normal"

# The template sheet only had formatting defined for the index column (A)
# down to row 3; stamp the same format across the full A2:A13 range so
# every row in column A is styled consistently (matches the other
# quarterly sheets' index-column formatting).
$ws.Range("A2").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at the
#    top of the data and renumber the existing index column.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the index column's formatting (copied from the row below, which
# still carries the original style) before writing the new row's values.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)
$summary.Application.CutCopyMode = $false

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 1.21

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 3) Restore "2020-Q4" as the selected tab (it was the selected tab before
#    this edit, and inserting/copying sheets shifts focus to the newest one).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("总计").Activate()
